$wb = $excel.ActiveWorkbook

# Update both "展览" and "全部类型" sheets which contain the same data table.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8711
    $ws.Range("F4").Value = 412
    $ws.Range("F5").Value = 46
}
